# Add sleep in the Tests that have failed due to SleepTiming
#
# This script:
#  1) Renames the "Sanity Suite Test" family of test-data strings to the
#     "Web Data 3" family everywhere they are used in the workbook.
#  2) Updates the recorded test-cursor cell selection on several sheets
#     (reflecting where execution stopped/paused after adding sleeps).
#  3) Updates one data cell (Mediation!T1) to the "New Test Code
#     Description" value.
#  4) Leaves the Mediation sheet active/selected (last-touched sheet),
#     matching the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Global text replacements (applied across every sheet/cell so the
#    shared string table is fully renamed, not just partially aliased).
# ---------------------------------------------------------------------
$replacements = @{
    "Sanity Suite Test" = "Web Data 3";
    "Sanity Suite Test Child" = "Web Data 3 Child";
    "Successfully created Sanity Suite Test Child. You can now login with the username admin after your password is set. Password reset link is sent to your email." = "Successfully created Web Data 3 Child. You can now login with the username admin after your password is set. Password reset link is sent to your email.";
    "Sanity Suite Test Reseller2" = "Web Data 3 Reseller";
    "Successfully created Sanity Suite Test Reseller2. You can now login with the username admin after your password is set. Password reset link is sent to your email." = "Successfully created Web Data 3 Reseller. You can now login with the username admin after your password is set. Password reset link is sent to your email.";
    "Working as admin Sanity Suite Test Child X" = "Working as admin Web Data 3 Child X";
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ($replacements.ContainsKey($v)) {
            $cell.Value = $replacements[$v]
        }
    }
}

# ---------------------------------------------------------------------
# 2) Update the recorded test-cursor cell / selection on each sheet.
#    (Mediation is done last so it ends up the active/selected tab.)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LoginData")
[void]$ws.Range("C1").Select()

$ws = $wb.Worksheets.Item("SignupChildComp")
[void]$ws.Range("S8").Select()

$ws = $wb.Worksheets.Item("SignupChildCompInvoice")
[void]$ws.Range("U1").Select()

$ws = $wb.Worksheets.Item("ConfigCollection")
[void]$ws.Range("C1").Select()

$ws = $wb.Worksheets.Item("ImpersonateInfo")
[void]$ws.Range("F4").Select()

$ws = $wb.Worksheets.Item("AddProduct")
[void]$ws.Range("J1").Select()

# ---------------------------------------------------------------------
# 3) Mediation sheet: update T1's value and finish with it selected /
#    active, matching the workbook's new activeTab.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mediation")
$ws.Activate()
$ws.Range("T1").Value = "New Test Code Description"
[void]$ws.Range("T1").Select()
